$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-09 07:48:25"
$ws.Range("I2").Value = "0.1 mm"
$ws.Range("O2").Value = "-3.8 °C"
$ws.Range("E3").Value = "2026-02-09 07:48:28"
$ws.Range("I3").Value = "0.5 mm"
$ws.Range("O3").Value = "-5.9 °C"
$ws.Range("E4").Value = "2026-02-09 07:48:30"
$ws.Range("J4").Value = "1008.1 hPa"
$ws.Range("K4").Value = "0.0 MJ/m2"
$ws.Range("N4").Value = "2.6 °C 7:13 TU"
$ws.Range("O4").Value = "4.0 °C"
$ws.Range("E5").Value = "2026-02-09 07:48:33"
$ws.Range("E6").Value = "2026-02-09 07:48:35"
$ws.Range("N6").Value = "4.4 °C 7:14 TU"
$ws.Range("O6").Value = "6.0 °C"
$ws.Range("E7").Value = "2026-02-09 07:48:37"
$ws.Range("H7").Value = "67%"
$ws.Range("L7").Value = "15.5 km/h - 312º 7:17 TU"
$ws.Range("N7").Value = "10.5 °C 7:19 TU"
$ws.Range("O7").Value = "11.0 °C"
$ws.Range("E8").Value = "2026-02-09 07:48:40"
$ws.Range("L8").Value = "24.1 km/h - 298º 7:17 TU"
$ws.Range("N8").Value = "6.1 °C 7:01 TU"
$ws.Range("O8").Value = "7.0 °C"
$ws.Range("E9").Value = "2026-02-09 07:48:42"
$ws.Range("H9").Value = "91%"
$ws.Range("N9").Value = "1.8 °C 7:09 TU"
$ws.Range("O9").Value = "5.6 °C"
$ws.Range("E10").Value = "2026-02-09 07:48:45"
$ws.Range("K10").Value = "0.0 MJ/m2"
$ws.Range("O10").Value = "4.7 °C"
$ws.Range("E11").Value = "2026-02-09 07:48:47"
$ws.Range("H11").Value = "96%"
$ws.Range("E12").Value = "2026-02-09 07:48:50"
$ws.Range("H12").Value = "94%"
$ws.Range("O12").Value = "6.3 °C"
$ws.Range("E13").Value = "2026-02-09 07:48:52"
$ws.Range("K13").Value = "0.0 MJ/m2"
$ws.Range("E14").Value = "2026-02-09 07:48:55"
$ws.Range("H14").Value = "94%"
$ws.Range("E15").Value = "2026-02-09 07:48:57"
$ws.Range("N15").Value = "2.3 °C 7:09 TU"
$ws.Range("O15").Value = "4.7 °C"
$ws.Range("E16").Value = "2026-02-09 07:48:59"
$ws.Range("H16").Value = "67%"
$ws.Range("I16").Value = "0.3 mm"
$ws.Range("O16").Value = "-5.5 °C"
$ws.Range("E17").Value = "2026-02-09 07:49:02"
$ws.Range("K17").Value = "0.1 MJ/m2"
$ws.Range("L17").Value = "36.7 km/h - 258º 7:16 TU"
$ws.Range("N17").Value = "-1.4 °C 6:50 TU"
$ws.Range("O17").Value = "-0.3 °C"
$ws.Range("E18").Value = "2026-02-09 07:49:04"
$ws.Range("J18").Value = "1008.0 hPa"
$ws.Range("O18").Value = "5.5 °C"
$ws.Range("E19").Value = "2026-02-09 07:49:07"
$ws.Range("N19").Value = "2.1 °C 7:00 TU"
$ws.Range("O19").Value = "3.0 °C"
$ws.Range("E20").Value = "2026-02-09 07:49:09"
$ws.Range("E21").Value = "2026-02-09 07:49:11"
$ws.Range("J21").Value = "1010.0 hPa"
$ws.Range("K21").Value = "0.0 MJ/m2"
$ws.Range("O21").Value = "0.3 °C"
$ws.Range("E22").Value = "2026-02-09 07:49:14"
$ws.Range("E23").Value = "2026-02-09 07:49:16"
$ws.Range("E24").Value = "2026-02-09 07:49:18"
$ws.Range("H24").Value = "91%"
$ws.Range("E25").Value = "2026-02-09 07:49:21"
$ws.Range("H25").Value = "73%"
$ws.Range("K25").Value = "0.0 MJ/m2"
$ws.Range("O25").Value = "-4.0 °C"
$ws.Range("E26").Value = "2026-02-09 07:49:23"
$ws.Range("J26").Value = "1008.5 hPa"
$ws.Range("K26").Value = "0.0 MJ/m2"
$ws.Range("E27").Value = "2026-02-09 07:49:26"
$ws.Range("K27").Value = "0.0 MJ/m2"
$ws.Range("E28").Value = "2026-02-09 07:49:28"
$ws.Range("J28").Value = "1008.4 hPa"
$ws.Range("E29").Value = "2026-02-09 07:49:31"
$ws.Range("H29").Value = "97%"
$ws.Range("K29").Value = "0.0 MJ/m2"
$ws.Range("N29").Value = "2.2 °C 6:45 TU"
$ws.Range("O29").Value = "4.6 °C"
$ws.Range("E30").Value = "2026-02-09 07:49:33"
$ws.Range("O30").Value = "6.0 °C"
$ws.Range("E31").Value = "2026-02-09 07:49:36"
$ws.Range("K31").Value = "0.0 MJ/m2"
$ws.Range("E32").Value = "2026-02-09 07:49:38"
$ws.Range("K32").Value = "0.0 MJ/m2"
$ws.Range("E33").Value = "2026-02-09 07:49:41"
$ws.Range("O33").Value = "-0.9 °C"
$ws.Range("E34").Value = "2026-02-09 07:49:43"
$ws.Range("E35").Value = "2026-02-09 07:49:46"
$ws.Range("H35").Value = "70%"
$ws.Range("I35").Value = "0.2 mm"
$ws.Range("J35").Value = "1010.0 hPa"
$ws.Range("N35").Value = "1.9 °C 7:28 TU"
$ws.Range("O35").Value = "3.4 °C"
$ws.Range("E36").Value = "2026-02-09 07:49:48"
$ws.Range("H36").Value = "88%"
$ws.Range("J36").Value = "1007.6 hPa"
$ws.Range("O36").Value = "7.4 °C"
$ws.Range("E37").Value = "2026-02-09 07:49:50"
$ws.Range("J37").Value = "1009.4 hPa"
$ws.Range("L37").Value = "16.9 km/h - 241º 7:29 TU"
$ws.Range("E38").Value = "2026-02-09 07:49:53"
$ws.Range("I38").Value = "0.1 mm"
$ws.Range("K38").Value = "0.0 MJ/m2"
$ws.Range("O38").Value = "5.3 °C"
$ws.Range("E39").Value = "2026-02-09 07:49:55"
$ws.Range("K39").Value = "0.0 MJ/m2"
$ws.Range("E40").Value = "2026-02-09 07:49:58"
$ws.Range("E41").Value = "2026-02-09 07:50:00"
$ws.Range("E42").Value = "2026-02-09 07:50:03"
$ws.Range("N42").Value = "3.1 °C 7:01 TU"
$ws.Range("O42").Value = "5.6 °C"
$ws.Range("E43").Value = "2026-02-09 07:50:05"
$ws.Range("K43").Value = "0.1 MJ/m2"
$ws.Range("E44").Value = "2026-02-09 07:50:07"
$ws.Range("I44").Value = "0.1 mm"
$ws.Range("E45").Value = "2026-02-09 07:50:10"
$ws.Range("J45").Value = "1009.2 hPa"
$ws.Range("K45").Value = "0.0 MJ/m2"
$ws.Range("E46").Value = "2026-02-09 07:50:12"
